$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Crime complaints table updates (rows 14, 16-25, 27, 28, 29) ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = -50
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -30.76923076923
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 17
$ws.Range("K16").Value = -29.411764705882
$ws.Range("L16").Value = 9.090909090909
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -83.098591549295
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -53.333333333333
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = -48
$ws.Range("L17").Value = -63.888888888888
$ws.Range("M17").Value = -13.333333333333
$ws.Range("N17").Value = -40.90909090909
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 5
$ws.Range("L18").Value = -16.666666666666
$ws.Range("M18").Value = -85.714285714285
$ws.Range("N18").Value = -94.117647058823
$ws.Range("C19").Value = 11
$ws.Range("E19").Value = 120
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 19
$ws.Range("H19").Value = 21.052631578947
$ws.Range("I19").Value = 37
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 23.333333333333
$ws.Range("L19").Value = 85
$ws.Range("M19").Value = 32.142857142857
$ws.Range("N19").Value = -17.777777777777
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 27.272727272727
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 53.846153846153
$ws.Range("L20").Value = 11.111111111111
$ws.Range("M20").Value = 11.111111111111
$ws.Range("N20").Value = -92.647058823529
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 88
$ws.Range("J21").Value = 92
$ws.Range("K21").Value = -4.347826086956
$ws.Range("L21").Value = -7.368421052631
$ws.Range("M21").Value = -24.137931034482
$ws.Range("N21").Value = -82.293762575452
$ws.Range("F22").Value = 3
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 100
$ws.Range("D23").Value = 2
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 20
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = 44.444444444444
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -30
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 53
$ws.Range("H24").Value = -22.641509433962
$ws.Range("I24").Value = 63
$ws.Range("J24").Value = 79
$ws.Range("K24").Value = -20.253164556962
$ws.Range("L24").Value = 5
$ws.Range("M24").Value = 40
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -3.846153846153
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = 8.333333333333
$ws.Range("L25").Value = 11.428571428571
$ws.Range("M25").Value = -15.217391304347
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("C26").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 2
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -75
$ws.Range("N28").Value = -50
$ws.Range("C26").Copy($ws.Range("C29"))
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -60
$ws.Range("N29").Value = -50
